$d = $word.ActiveDocument

# --- Paragraph 1: update the ID placeholder text and drop the trailing run ---
$d.Content.Find.Execute("**ID__AFFARS_5336_topic_14__ID**", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_5336_602_3__ID**", 2)

$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("**ID__AFFARS_5336_602_3__ID** ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "**ID__AFFARS_5336_602_3__ID**", 2)

# --- Paragraph 1: add a paragraph border (5 twip space on every side) ---
$b = $p1.ParagraphFormat.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# --- Paragraph 1: widen the left indent from 120 to 225 twips (6pt -> 11.25pt) ---
$p1.ParagraphFormat.LeftIndent = 11.25
